$d = $word.ActiveDocument

# Locate the "3. Microservicios Backend" heading paragraph
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "3. Microservicios Backend") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Heading paragraph '3. Microservicios Backend' not found"
}

# Insert a new paragraph right after the heading
$target.Range.InsertParagraphAfter()

# Re-find the paragraph that was just created (the one following the heading)
$newPara = $target.Next()
$r = $newPara.Range

# Use the default ("Normal") paragraph style, as in the target markup
# (no explicit w:pStyle on the inserted paragraph).
$r.Style = "Normal"

$line1 = "▪ api-bff: servicio Backend for Frontend que actúa como orquestador para el Dashboard General de Bienvenida. "
$line2 = "   Centraliza y unifica las llamadas a múltiples microservicios (api-cuidados, api-citas, api-rutinas, profile-service) "
$line3 = "   proporcionando un único endpoint optimizado que devuelve rutinas del día, próximas citas, últimos cuidados y estadísticas rápidas."
$line4 = "   Mejora el rendimiento reduciendo la latencia y simplificando la lógica en el frontend."

# [char]11 is a manual line break (w:br) when assigned through Range.Text
$r.Text = $line1 + [char]11 + $line2 + [char]11 + $line3 + [char]11 + $line4
